$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I27").Value = 0.0195821
$ws.Range("J27").Value = 0.0313632
$ws.Range("I28").Value = 0.0370894
$ws.Range("J28").Value = 0.0146784
